# Update of 2025 data and RF changes
# Column I ("RF") for rows 21-44 changes from 13.21378378378378 to 15.73978723404255

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I21:I44").Value = 15.73978723404255
